$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update existing rows 8, 9, 12, 13, 14, 15 with revised field values
# ---------------------------------------------------------------------------

# Row 8
$ws.Cells.Item(8, 4).Value = 44434            # Fecha
$ws.Cells.Item(8, 13).Value = 40              # Volumen
$ws.Cells.Item(8, 14).Value = 35000           # Precio minimo
$ws.Cells.Item(8, 15).Value = 35000           # Precio maximo
$ws.Cells.Item(8, 16).Value = 35000           # Precio promedio ponderado
$ws.Cells.Item(8, 17).Value = "$/caja 18 kilos"  # Unidad de comercializacion
$ws.Cells.Item(8, 19).Value = 1944            # Precio $/Kg
$ws.Cells.Item(8, 20).Value = 18              # Kg / unidad

# Row 9
$ws.Cells.Item(9, 4).Value = 44363            # Fecha
$ws.Cells.Item(9, 13).Value = 144             # Volumen
$ws.Cells.Item(9, 14).Value = 1700            # Precio minimo
$ws.Cells.Item(9, 15).Value = 1700            # Precio maximo
$ws.Cells.Item(9, 16).Value = 1700            # Precio promedio ponderado
$ws.Cells.Item(9, 18).Value = "Región de Arica y Parinacota"  # Origen
$ws.Cells.Item(9, 19).Value = 1700            # Precio $/Kg

# Row 12
$ws.Cells.Item(12, 4).Value = 44435           # Fecha
$ws.Cells.Item(12, 18).Value = "Perú"         # Origen

# Row 13
$ws.Cells.Item(13, 4).Value = 44435           # Fecha
$ws.Cells.Item(13, 13).Value = 105            # Volumen

# Row 14
$ws.Cells.Item(14, 4).Value = 44431           # Fecha
$ws.Cells.Item(14, 13).Value = 30             # Volumen
$ws.Cells.Item(14, 18).Value = "Región de Arica y Parinacota"  # Origen

# Row 15
$ws.Cells.Item(15, 4).Value = 44405           # Fecha
$ws.Cells.Item(15, 13).Value = 10             # Volumen

# ---------------------------------------------------------------------------
# Append new rows 16-20, cloning row 15's layout/format then overwriting
# the varying fields (Fecha, Volumen, Origen)
# ---------------------------------------------------------------------------

$ws.Range("A15:T15").Copy($ws.Range("A16:T16"))
$ws.Cells.Item(16, 4).Value = 44294
$ws.Cells.Item(16, 13).Value = 15
$ws.Cells.Item(16, 18).Value = "Región de Arica y Parinacota"

$ws.Range("A15:T15").Copy($ws.Range("A17:T17"))
$ws.Cells.Item(17, 4).Value = 44369
$ws.Cells.Item(17, 13).Value = 5
$ws.Cells.Item(17, 18).Value = "Perú"

$ws.Range("A15:T15").Copy($ws.Range("A18:T18"))
$ws.Cells.Item(18, 4).Value = 44433
$ws.Cells.Item(18, 13).Value = 15
$ws.Cells.Item(18, 18).Value = "Región de Arica y Parinacota"

$ws.Range("A15:T15").Copy($ws.Range("A19:T19"))
$ws.Cells.Item(19, 4).Value = 44418
$ws.Cells.Item(19, 13).Value = 30
$ws.Cells.Item(19, 18).Value = "Región de Arica y Parinacota"

$ws.Range("A15:T15").Copy($ws.Range("A20:T20"))
$ws.Cells.Item(20, 4).Value = 44432
$ws.Cells.Item(20, 13).Value = 10
$ws.Cells.Item(20, 18).Value = "Perú"
